$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 78/79 and 80/81 (all columns except A) ---
# Row 78
$ws.Cells.Item(78, 2).Value = 5499423
$ws.Cells.Item(78, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(78, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(78, 5).Value = 45066.41666666666
$ws.Cells.Item(78, 6).Value = "Olimpija Ljubljana"
$ws.Cells.Item(78, 7).Value = "NK Celje"
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 2
$ws.Cells.Item(78, 10).Value = "A"
$ws.Cells.Item(78, 11).Value = 2.5
$ws.Cells.Item(78, 12).Value = 3.3
$ws.Cells.Item(78, 13).Value = 2.5
$ws.Cells.Item(78, 14).Value = 2.55
$ws.Cells.Item(78, 15).Value = 3.25
$ws.Cells.Item(78, 16).Value = 2.45
$ws.Cells.Item(78, 17).Value = 0
$ws.Cells.Item(78, 18).Value = 1.95
$ws.Cells.Item(78, 19).Value = 1.85
$ws.Cells.Item(78, 20).Value = 2.5
$ws.Cells.Item(78, 21).Value = 1.85
$ws.Cells.Item(78, 22).Value = 1.95
$ws.Cells.Item(78, 23).Value = -1
$ws.Cells.Item(78, 24).Value = -1
$ws.Cells.Item(78, 25).Value = 1.45
$ws.Cells.Item(78, 26).Value = -1
$ws.Cells.Item(78, 27).Value = 0.8500000000000001
$ws.Cells.Item(78, 28).Value = -1
$ws.Cells.Item(78, 29).Value = 0.95

# Row 79
$ws.Cells.Item(79, 2).Value = 5498504
$ws.Cells.Item(79, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(79, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(79, 5).Value = 45066.41666666666
$ws.Cells.Item(79, 6).Value = "NK Maribor"
$ws.Cells.Item(79, 7).Value = "NK Bravo"
$ws.Cells.Item(79, 8).Value = 1
$ws.Cells.Item(79, 9).Value = 1
$ws.Cells.Item(79, 10).Value = "D"
$ws.Cells.Item(79, 11).Value = 1.571
$ws.Cells.Item(79, 12).Value = 3.8
$ws.Cells.Item(79, 13).Value = 4.75
$ws.Cells.Item(79, 14).Value = 1.533
$ws.Cells.Item(79, 15).Value = 4
$ws.Cells.Item(79, 16).Value = 4.75
$ws.Cells.Item(79, 17).Value = -1
$ws.Cells.Item(79, 18).Value = 1.95
$ws.Cells.Item(79, 19).Value = 1.85
$ws.Cells.Item(79, 20).Value = 2.75
$ws.Cells.Item(79, 21).Value = 1.875
$ws.Cells.Item(79, 22).Value = 1.925
$ws.Cells.Item(79, 23).Value = -1
$ws.Cells.Item(79, 24).Value = 3
$ws.Cells.Item(79, 25).Value = -1
$ws.Cells.Item(79, 26).Value = -1
$ws.Cells.Item(79, 27).Value = 0.8500000000000001
$ws.Cells.Item(79, 28).Value = -1
$ws.Cells.Item(79, 29).Value = 0.925

# Row 80
$ws.Cells.Item(80, 2).Value = 5495053
$ws.Cells.Item(80, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(80, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(80, 5).Value = 45066.63541666666
$ws.Cells.Item(80, 6).Value = "NK Radomlje"
$ws.Cells.Item(80, 7).Value = "NK Domzale"
$ws.Cells.Item(80, 8).Value = 1
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = "H"
$ws.Cells.Item(80, 11).Value = 2.55
$ws.Cells.Item(80, 12).Value = 3.1
$ws.Cells.Item(80, 13).Value = 2.55
$ws.Cells.Item(80, 14).Value = 3.75
$ws.Cells.Item(80, 15).Value = 3.4
$ws.Cells.Item(80, 16).Value = 1.833
$ws.Cells.Item(80, 17).Value = 0.5
$ws.Cells.Item(80, 18).Value = 1.925
$ws.Cells.Item(80, 19).Value = 1.875
$ws.Cells.Item(80, 20).Value = 2.5
$ws.Cells.Item(80, 21).Value = 1.975
$ws.Cells.Item(80, 22).Value = 1.825
$ws.Cells.Item(80, 23).Value = 2.75
$ws.Cells.Item(80, 24).Value = -1
$ws.Cells.Item(80, 25).Value = -1
$ws.Cells.Item(80, 26).Value = 0.925
$ws.Cells.Item(80, 27).Value = -1
$ws.Cells.Item(80, 28).Value = -1
$ws.Cells.Item(80, 29).Value = 0.825

# Row 81
$ws.Cells.Item(81, 2).Value = 5498503
$ws.Cells.Item(81, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(81, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(81, 5).Value = 45066.63541666666
$ws.Cells.Item(81, 6).Value = "FC Koper"
$ws.Cells.Item(81, 7).Value = "NS Mura"
$ws.Cells.Item(81, 8).Value = 1
$ws.Cells.Item(81, 9).Value = 2
$ws.Cells.Item(81, 10).Value = "A"
$ws.Cells.Item(81, 11).Value = 2.05
$ws.Cells.Item(81, 12).Value = 3.3
$ws.Cells.Item(81, 13).Value = 3.25
$ws.Cells.Item(81, 14).Value = 2
$ws.Cells.Item(81, 15).Value = 3.4
$ws.Cells.Item(81, 16).Value = 3.25
$ws.Cells.Item(81, 17).Value = -0.5
$ws.Cells.Item(81, 18).Value = 2
$ws.Cells.Item(81, 19).Value = 1.8
$ws.Cells.Item(81, 20).Value = 2.5
$ws.Cells.Item(81, 21).Value = 1.825
$ws.Cells.Item(81, 22).Value = 1.975
$ws.Cells.Item(81, 23).Value = -1
$ws.Cells.Item(81, 24).Value = -1
$ws.Cells.Item(81, 25).Value = 2.25
$ws.Cells.Item(81, 26).Value = -1
$ws.Cells.Item(81, 27).Value = 0.8
$ws.Cells.Item(81, 28).Value = 0.825
$ws.Cells.Item(81, 29).Value = -1

# --- Copy A/E column formatting (style) down into the new rows 183-188 ---
$ws.Range("A180").Copy()
$ws.Range("A183").PasteSpecial(-4122)
$ws.Range("A184").PasteSpecial(-4122)
$ws.Range("A185").PasteSpecial(-4122)
$ws.Range("A186").PasteSpecial(-4122)
$ws.Range("A187").PasteSpecial(-4122)
$ws.Range("A188").PasteSpecial(-4122)
$ws.Range("E180").Copy()
$ws.Range("E183").PasteSpecial(-4122)
$ws.Range("E184").PasteSpecial(-4122)
$ws.Range("E185").PasteSpecial(-4122)
$ws.Range("E186").PasteSpecial(-4122)
$ws.Range("E187").PasteSpecial(-4122)
$ws.Range("E188").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update rows 181, 182 (new B-AC data incl. added H/I/J score) and new rows 183-188 ---
# Row 181
$ws.Cells.Item(181, 1).Value = 179
$ws.Cells.Item(181, 2).Value = 7680775
$ws.Cells.Item(181, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(181, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(181, 5).Value = 45339.375
$ws.Cells.Item(181, 6).Value = "NK Rogaska"
$ws.Cells.Item(181, 7).Value = "NK Domzale"
$ws.Cells.Item(181, 8).Value = 4
$ws.Cells.Item(181, 9).Value = 1
$ws.Cells.Item(181, 10).Value = "H"
$ws.Cells.Item(181, 11).Value = 3.5
$ws.Cells.Item(181, 12).Value = 3.4
$ws.Cells.Item(181, 13).Value = 1.95
$ws.Cells.Item(181, 14).Value = 3.3
$ws.Cells.Item(181, 15).Value = 3.2
$ws.Cells.Item(181, 16).Value = 2.1
$ws.Cells.Item(181, 17).Value = 0.25
$ws.Cells.Item(181, 18).Value = 1.975
$ws.Cells.Item(181, 19).Value = 1.825
$ws.Cells.Item(181, 20).Value = 2.75
$ws.Cells.Item(181, 21).Value = 2.05
$ws.Cells.Item(181, 22).Value = 1.75
$ws.Cells.Item(181, 23).Value = 2.3
$ws.Cells.Item(181, 24).Value = -1
$ws.Cells.Item(181, 25).Value = -1
$ws.Cells.Item(181, 26).Value = 0.9750000000000001
$ws.Cells.Item(181, 27).Value = -1
$ws.Cells.Item(181, 28).Value = 1.05
$ws.Cells.Item(181, 29).Value = -1

# Row 182
$ws.Cells.Item(182, 1).Value = 180
$ws.Cells.Item(182, 2).Value = 7680774
$ws.Cells.Item(182, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(182, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(182, 5).Value = 45339.45833333334
$ws.Cells.Item(182, 6).Value = "NK Radomlje"
$ws.Cells.Item(182, 7).Value = "Olimpija Ljubljana"
$ws.Cells.Item(182, 8).Value = 1
$ws.Cells.Item(182, 9).Value = 3
$ws.Cells.Item(182, 10).Value = "A"
$ws.Cells.Item(182, 11).Value = 4.8
$ws.Cells.Item(182, 12).Value = 3.8
$ws.Cells.Item(182, 13).Value = 1.615
$ws.Cells.Item(182, 14).Value = 5.25
$ws.Cells.Item(182, 15).Value = 3.75
$ws.Cells.Item(182, 16).Value = 1.571
$ws.Cells.Item(182, 17).Value = 0.75
$ws.Cells.Item(182, 18).Value = 2.025
$ws.Cells.Item(182, 19).Value = 1.775
$ws.Cells.Item(182, 20).Value = 2.5
$ws.Cells.Item(182, 21).Value = 1.85
$ws.Cells.Item(182, 22).Value = 1.95
$ws.Cells.Item(182, 23).Value = -1
$ws.Cells.Item(182, 24).Value = -1
$ws.Cells.Item(182, 25).Value = 0.571
$ws.Cells.Item(182, 26).Value = -1
$ws.Cells.Item(182, 27).Value = 0.7749999999999999
$ws.Cells.Item(182, 28).Value = 0.8500000000000001
$ws.Cells.Item(182, 29).Value = -1

# Row 183
$ws.Cells.Item(183, 1).Value = 181
$ws.Cells.Item(183, 2).Value = 7680773
$ws.Cells.Item(183, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(183, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(183, 5).Value = 45340.375
$ws.Cells.Item(183, 6).Value = "NK Bravo"
$ws.Cells.Item(183, 7).Value = "NK Celje"
$ws.Cells.Item(183, 8).Value = 0
$ws.Cells.Item(183, 9).Value = 0
$ws.Cells.Item(183, 10).Value = "D"
$ws.Cells.Item(183, 11).Value = 4.5
$ws.Cells.Item(183, 12).Value = 3.4
$ws.Cells.Item(183, 13).Value = 1.75
$ws.Cells.Item(183, 14).Value = 5
$ws.Cells.Item(183, 15).Value = 3.6
$ws.Cells.Item(183, 16).Value = 1.7
$ws.Cells.Item(183, 17).Value = 0.75
$ws.Cells.Item(183, 18).Value = 1.925
$ws.Cells.Item(183, 19).Value = 1.875
$ws.Cells.Item(183, 20).Value = 2.5
$ws.Cells.Item(183, 21).Value = 1.95
$ws.Cells.Item(183, 22).Value = 1.85
$ws.Cells.Item(183, 23).Value = -1
$ws.Cells.Item(183, 24).Value = 2.6
$ws.Cells.Item(183, 25).Value = -1
$ws.Cells.Item(183, 26).Value = 0.925
$ws.Cells.Item(183, 27).Value = -1
$ws.Cells.Item(183, 28).Value = -1
$ws.Cells.Item(183, 29).Value = 0.8500000000000001

# Row 184
$ws.Cells.Item(184, 1).Value = 182
$ws.Cells.Item(184, 2).Value = 6814411
$ws.Cells.Item(184, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(184, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(184, 5).Value = 45342.5625
$ws.Cells.Item(184, 6).Value = "NK Radomlje"
$ws.Cells.Item(184, 7).Value = "NK Aluminij"
$ws.Cells.Item(184, 11).Value = 1.909
$ws.Cells.Item(184, 12).Value = 3.4
$ws.Cells.Item(184, 13).Value = 3.5
$ws.Cells.Item(184, 14).Value = 1.909
$ws.Cells.Item(184, 15).Value = 3.4
$ws.Cells.Item(184, 16).Value = 3.5
$ws.Cells.Item(184, 17).Value = -0.5
$ws.Cells.Item(184, 18).Value = 1.95
$ws.Cells.Item(184, 19).Value = 1.85
$ws.Cells.Item(184, 20).Value = 2.5
$ws.Cells.Item(184, 21).Value = 1.9
$ws.Cells.Item(184, 22).Value = 1.9
$ws.Cells.Item(184, 23).Value = 0
$ws.Cells.Item(184, 24).Value = 0
$ws.Cells.Item(184, 25).Value = 0
$ws.Cells.Item(184, 26).Value = 0
$ws.Cells.Item(184, 27).Value = 0

# Row 185
$ws.Cells.Item(185, 1).Value = 183
$ws.Cells.Item(185, 2).Value = 6814408
$ws.Cells.Item(185, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(185, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(185, 5).Value = 45343.45833333334
$ws.Cells.Item(185, 6).Value = "Olimpija Ljubljana"
$ws.Cells.Item(185, 7).Value = "NK Domzale"
$ws.Cells.Item(185, 11).Value = 1.533
$ws.Cells.Item(185, 12).Value = 4
$ws.Cells.Item(185, 13).Value = 5
$ws.Cells.Item(185, 14).Value = 1.4
$ws.Cells.Item(185, 15).Value = 4.5
$ws.Cells.Item(185, 16).Value = 6.5
$ws.Cells.Item(185, 17).Value = -1.25
$ws.Cells.Item(185, 18).Value = 1.925
$ws.Cells.Item(185, 19).Value = 1.875
$ws.Cells.Item(185, 20).Value = 2.75
$ws.Cells.Item(185, 21).Value = 1.8
$ws.Cells.Item(185, 22).Value = 2
$ws.Cells.Item(185, 23).Value = 0
$ws.Cells.Item(185, 24).Value = 0
$ws.Cells.Item(185, 25).Value = 0
$ws.Cells.Item(185, 26).Value = 0
$ws.Cells.Item(185, 27).Value = 0

# Row 186
$ws.Cells.Item(186, 1).Value = 184
$ws.Cells.Item(186, 2).Value = 6814410
$ws.Cells.Item(186, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(186, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(186, 5).Value = 45343.5625
$ws.Cells.Item(186, 6).Value = "FC Koper"
$ws.Cells.Item(186, 7).Value = "NS Mura"
$ws.Cells.Item(186, 11).Value = 1.8
$ws.Cells.Item(186, 12).Value = 3.4
$ws.Cells.Item(186, 13).Value = 4
$ws.Cells.Item(186, 14).Value = 1.8
$ws.Cells.Item(186, 15).Value = 3.4
$ws.Cells.Item(186, 16).Value = 4
$ws.Cells.Item(186, 17).Value = -0.5
$ws.Cells.Item(186, 18).Value = 1.825
$ws.Cells.Item(186, 19).Value = 1.975
$ws.Cells.Item(186, 20).Value = 2.5
$ws.Cells.Item(186, 21).Value = 1.95
$ws.Cells.Item(186, 22).Value = 1.85
$ws.Cells.Item(186, 23).Value = 0
$ws.Cells.Item(186, 24).Value = 0
$ws.Cells.Item(186, 25).Value = 0
$ws.Cells.Item(186, 26).Value = 0
$ws.Cells.Item(186, 27).Value = 0

# Row 187
$ws.Cells.Item(187, 1).Value = 185
$ws.Cells.Item(187, 2).Value = 6816452
$ws.Cells.Item(187, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(187, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(187, 5).Value = 45344.45833333334
$ws.Cells.Item(187, 6).Value = "NK Celje"
$ws.Cells.Item(187, 7).Value = "NK Rogaska"
$ws.Cells.Item(187, 11).Value = 1.25
$ws.Cells.Item(187, 12).Value = 5.75
$ws.Cells.Item(187, 13).Value = 8.5
$ws.Cells.Item(187, 14).Value = 1.25
$ws.Cells.Item(187, 15).Value = 5.75
$ws.Cells.Item(187, 16).Value = 8.5
$ws.Cells.Item(187, 17).Value = -1.75
$ws.Cells.Item(187, 18).Value = 1.9
$ws.Cells.Item(187, 19).Value = 1.9
$ws.Cells.Item(187, 20).Value = 3.25
$ws.Cells.Item(187, 21).Value = 1.975
$ws.Cells.Item(187, 22).Value = 1.825
$ws.Cells.Item(187, 23).Value = 0
$ws.Cells.Item(187, 24).Value = 0
$ws.Cells.Item(187, 25).Value = 0
$ws.Cells.Item(187, 26).Value = 0
$ws.Cells.Item(187, 27).Value = 0

# Row 188
$ws.Cells.Item(188, 1).Value = 186
$ws.Cells.Item(188, 2).Value = 6814409
$ws.Cells.Item(188, 3).Value = "Slovenia Prva Liga"
$ws.Cells.Item(188, 4).Value = "Slovenia Prva Liga"
$ws.Cells.Item(188, 5).Value = 45344.5625
$ws.Cells.Item(188, 6).Value = "NK Maribor"
$ws.Cells.Item(188, 7).Value = "NK Bravo"
$ws.Cells.Item(188, 11).Value = 1.615
$ws.Cells.Item(188, 12).Value = 3.75
$ws.Cells.Item(188, 13).Value = 4.75
$ws.Cells.Item(188, 14).Value = 1.615
$ws.Cells.Item(188, 15).Value = 3.75
$ws.Cells.Item(188, 16).Value = 4.75
$ws.Cells.Item(188, 17).Value = -0.75
$ws.Cells.Item(188, 18).Value = 1.8
$ws.Cells.Item(188, 19).Value = 2
$ws.Cells.Item(188, 20).Value = 2.5
$ws.Cells.Item(188, 21).Value = 1.875
$ws.Cells.Item(188, 22).Value = 1.925
$ws.Cells.Item(188, 23).Value = 0
$ws.Cells.Item(188, 24).Value = 0
$ws.Cells.Item(188, 25).Value = 0
$ws.Cells.Item(188, 26).Value = 0
$ws.Cells.Item(188, 27).Value = 0

